# Add emerging technologies candidates to the "2020-S2" radar sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020-S2")

# Unicode "left single quotation mark" used (verbatim, matching source) in the
# PWA bullet about the phone's home screen.
$smartQuote = [char]0x2018

$pwaDescription = "<p>PWA features allow to close the gap to native applications and create similar user experiences. Features include</p>`n<ul><li>Work offline</li>`n<li>High performance</li>`n<li>Background processing in service workers in a separate thread</li>`n<li>Access to the phone's sensors</li>`n<li>Support for push notifications</li>`n<li>An icon on the phone${smartQuote}s home screen</li></ul>"

$amplifyDescription = "A set of products and tools that enables mobile and front-end web developers to build and deploy secure, scalable full stack applications, powered by AWS. With Amplify, you can configure app backends in minutes, connect them to your app in just a few lines of code, and deploy static web apps in three steps."

$gcpDescription = "A suite of cloud computing services that runs on the same infrastructure that Google uses internally for its end-user products, such as Google Search, Gmail, file storage, and YouTube."

$azureDevopsDescription = "A Software as a service (SaaS) platform from Microsoft that provides an end-to-end DevOps toolchain for developing and deploying software. It also integrates with most leading tools on the market and is a great option for orchestrating a DevOps toolchain."

$blazorDescription = "Blazor is a feature of ASP.NET for building interactive web UIs using C# instead of JavaScript. It's real .NET running in the browser on WebAssembly."

# --- Column A (name) first, in row order, so new shared strings are created
#     in the same order the source workbook has them (indices 71-74). ---
$ws.Cells.Item(16, 1).Value = "Progressive Web Apps"
$ws.Cells.Item(17, 1).Value = "AWS Amplify"
$ws.Cells.Item(18, 1).Value = "Google Cloud Platform"
$ws.Cells.Item(19, 1).Value = "Azure Devops"
$ws.Cells.Item(20, 1).Value = "Blazor"

# --- Columns B/C/D: all re-use shared strings that already exist in the
#     workbook, so the order here does not create new shared-string entries. ---
$ws.Cells.Item(16, 2).Value = "trial"
$ws.Cells.Item(16, 3).Value = "techniques"
$ws.Cells.Item(16, 4).Value = $true

$ws.Cells.Item(17, 2).Value = "trial"
$ws.Cells.Item(17, 3).Value = "languages & frameworks"
$ws.Cells.Item(17, 4).Value = $true

$ws.Cells.Item(18, 2).Value = "adopt"
$ws.Cells.Item(18, 3).Value = "platforms"
$ws.Cells.Item(18, 4).Value = $false

$ws.Cells.Item(19, 2).Value = "adopt"
$ws.Cells.Item(19, 3).Value = "platforms"
$ws.Cells.Item(19, 4).Value = $false

$ws.Cells.Item(20, 2).Value = "assess"
$ws.Cells.Item(20, 3).Value = "languages & frameworks"
$ws.Cells.Item(20, 4).Value = $true

# --- Column E (description) last, in row order, so the new shared strings
#     land at indices 75-79. ---
$ws.Cells.Item(16, 5).Value = $pwaDescription
$ws.Cells.Item(17, 5).Value = $amplifyDescription
$ws.Cells.Item(18, 5).Value = $gcpDescription
$ws.Cells.Item(19, 5).Value = $azureDevopsDescription
$ws.Cells.Item(20, 5).Value = $blazorDescription

# --- Formatting: rows 16-19 mirror the sheet's usual "new row" style
#     (A:D vertically centred + wrapped, E top-aligned + wrapped). Row 20's
#     A:D keep the worksheet's default (unstyled) formatting, matching the
#     other "assess" rows that were typed without the helper formatting. ---
$range16to19 = $ws.Range($ws.Cells.Item(16, 1), $ws.Cells.Item(19, 4))
$range16to19.WrapText = $true
$range16to19.VerticalAlignment = -4108

$descriptionRange = $ws.Range($ws.Cells.Item(16, 5), $ws.Cells.Item(20, 5))
$descriptionRange.WrapText = $true
$descriptionRange.VerticalAlignment = -4160

# The PWA description (row 16) contains embedded line breaks, which makes the
# headless engine stamp a (wrong, font-metric-less) customHeight on that row.
# Let the engine re-fit it so the row falls back to the sheet's normal
# (non-custom) row height, consistent with the other new rows.
$ws.Rows.Item(16).AutoFit()

[void]$ws.Range("E20").Select()
